# Fill in the hour log entries for rows 19-23 on the "DTT Test Hour Log" sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DTT Test Hour Log")

# Row 19
$ws.Range("A19").Value = "Fixed error regarding formCpomponent and changed components names"
$ws.Range("B19").Value = 1
$ws.Range("C19").Value = "12/13/2023"
$ws.Range("D19").Value = "Fixed error that meant tht editing was disabled, changed xomponents names and removed unused variables"

# Row 20
$ws.Range("A20").Value = "Fixing erorr with formComponent and updated the API part of the site"
$ws.Range("B20").Value = 4
$ws.Range("C20").Value = "12/18/2023"
$ws.Range("D20").Value = "Fixed error of editing part that came back, also removed the data.vue and added API.js"

# Row 21
$ws.Range("A21").Value = "Changed all the API calls from using the methods in vue to use the js file"
$ws.Range("B21").Value = 2
$ws.Range("C21").Value = "12/19/2023"
$ws.Range("D21").Value = "Changed API calls to use the js file instead of the vue file"

# Row 22
$ws.Range("A22").Value = "Fixed error with showing houses as a card and in detail form"
$ws.Range("B22").Value = 1
$ws.Range("C22").Value = "12/20/2023"
$ws.Range("D22").Value = "Changing the way empty housenumberAdditions are treated inside the website"

# Row 23
$ws.Range("A23").Value = "Fixing styling errors that came up by changing the way that the websites worked"
$ws.Range("B23").Value = 1
$ws.Range("C23").Value = "12/21/2023"
$ws.Range("D23").Value = "Fixed the bugs that showed on the screen due to it being different then first, Including formComponent and other files"

# Update the selected cell to match the saved view state
$ws.Range("D24").Select()
